$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.406.41'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.79%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.668.31'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.47%  '
$ws.Range('E4').Value = '  +0.43%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '219.96'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.05%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5248'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.59%  '
$ws.Range('E7').Value = '  +0.39%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2664'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.22%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06351'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.28%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.62'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.62%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07783'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.33%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.684.54'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.62%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.458'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.91%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5515'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.36%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0₅8239'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.32%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.44'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.30%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.428.31'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.97%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.004'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.20%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.730'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.62%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '193.14'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.03%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.29'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.08%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.256'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.17%  '
$ws.Range('E23').Value = '  +0.51%  '
$ws.Range('B24').Value = 'Monero'
$ws.Range('C24').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '138.74'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.90%  '
$ws.Range('B25').Value = 'Stellar'
$ws.Range('C25').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1257'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.28%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.372'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.09%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.18'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.56%  '
$ws.Range('E28').Value = '  +1.09%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.06110'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.95%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.291'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.78%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.593'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +5.30%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.386'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.45%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.678'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.36%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.000'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.02%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.422'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.23%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6037'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +7.52%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.768'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.93%  '
$ws.Range('E38').Value = '  +0.59%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.028'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.16%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.086.62'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +6.07%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8580'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.48%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.003'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.18%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '100.52'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.78%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.812.92'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.34%  '
$ws.Range('B45').Value = 'BabyDogeCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0₈110'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.27%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '57.85'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.11%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.150'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.57%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.005'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.28%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05203'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.26%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.483'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +8.08%  '
